$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.405.84"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.022.47"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.33"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.09"
$ws.Range("E8").Value = "  -6.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.61"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.319.78"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.814"
$ws.Range("E14").Value = "  -4.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.82"
$ws.Range("E15").Value = "  -9.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.040.92"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.326.71"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.97"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0848"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.23"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.52"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.64"
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.55"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.75"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -9.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0661"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.62"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.40"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.23"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("E40").Value = "  +3.93%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.21"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0960"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.413.03"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.20"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.99"
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.35"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.03"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.212.16"
$ws.Range("E51").Value = "  +0.99%  "
